# Update the header row of the "locations" sheet with nicer, human-readable
# column titles (id -> ID, store_name -> Store Name, store_address -> Store
# Address, store_postcode -> Store Postcode, kms -> Kilometers,
# tail_lift -> Tail Lift).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Store Name"
$ws.Range("C1").Value = "Store Address"
$ws.Range("D1").Value = "Store Postcode"
$ws.Range("E1").Value = "Kilometers"
$ws.Range("F1").Value = "Tail Lift"
